$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update name / email (status stays "ENVIADO") ---
$ws.Range("A2").Value = "Gleyse Oliveira Rosal"
$ws.Range("B2").Value = "gleyseor@gmail.com"

# --- Row 3: replace person, mark as sent, drop the "invalid e-mail" hyperlink ---
$ws.Range("B3").Hyperlinks.Delete()
$ws.Range("A3").Value = "Erasmo Rosa Dos Reis"
$ws.Range("B3").Value = "erasmo.reis@bol.com.br"
$ws.Range("C3").Style = "Normal"
$ws.Range("C3").Value = "ENVIADO"

# B3, and the left-over Email placeholders in B5/B6/B8, keep the "hyperlink
# color" look but without the underline now that the hyperlink itself is
# gone - copy that look from B4 (which already is no-underline) BEFORE B4
# itself gets its own formatting changed below.
$ws.Range("B4").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# --- Rows 5,6,8: clear the Name/Status placeholders, keep a plain (no-underline)
#     look on the Email placeholder cell ---
$ws.Range("A5").Clear()
$ws.Range("C5").Clear()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("A6").Clear()
$ws.Range("C6").Clear()
$ws.Range("B6").PasteSpecial(-4122)

# --- Row 7: fully cleared (disappears from the sheet) ---
$ws.Range("A7:C7").Clear()

$ws.Range("A8").Clear()
$ws.Range("C8").Clear()
$ws.Range("B8").PasteSpecial(-4122)

# --- Row 4: turn wrap text on for the placeholder cell ---
$ws.Range("B4").WrapText = $True
$ws.Range("B4").Font.Underline = $False

# --- Column A got a bit wider ---
$ws.Columns("A").ColumnWidth = 34.42

# --- Cursor ends up on A6 ---
$ws.Range("A6").Select()
